# Update the Flt3l-Flt3 NATMI output with the new TPM-derived values.
# Final table keeps only 3 data rows (ECs->MuSCs, FAPs->MuSCs, MuSCs->MuSCs);
# the former ECs->ECs, FAPs->ECs, MuSCs->ECs rows are dropped and the
# remaining rows are re-populated with the new TPM-based numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old rows 5-7 first so the sheet only spans A1:T7 -> A1:T4.
$ws.Range("A5:T7").EntireRow.Delete()

# ---- Row 2: ECs -> MuSCs -------------------------------------------------
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Flt3l"
$ws.Range("C2").Value = "Flt3"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.084097666666667
$ws.Range("H2").Value = 27.252293
$ws.Range("I2").Value = 0.1656252520634623
$ws.Range("J2").Value = 0.1656252520634623
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04013533333333334
$ws.Range("N2").Value = 0.120406
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.3645932878842222
$ws.Range("R2").Value = 3.281339590958
$ws.Range("S2").Value = 0.1656252520634623
$ws.Range("T2").Value = 0.1656252520634623

# ---- Row 3: FAPs -> MuSCs -------------------------------------------------
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Flt3l"
$ws.Range("C3").Value = "Flt3"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 40.138293
$ws.Range("H3").Value = 120.414879
$ws.Range("I3").Value = 0.7318189587410614
$ws.Range("J3").Value = 0.7318189587410614
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.04013533333333334
$ws.Range("N3").Value = 0.120406
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1.610963768986
$ws.Range("R3").Value = 14.498673920874
$ws.Range("S3").Value = 0.7318189587410614
$ws.Range("T3").Value = 0.7318189587410614

# ---- Row 4: MuSCs -> MuSCs -------------------------------------------------
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Flt3l"
$ws.Range("C4").Value = "Flt3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.624908
$ws.Range("H4").Value = 16.874724
$ws.Range("I4").Value = 0.1025557891954764
$ws.Range("J4").Value = 0.1025557891954764
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.04013533333333334
$ws.Range("N4").Value = 0.120406
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.2257575575493334
$ws.Range("R4").Value = 2.031818017944
$ws.Range("S4").Value = 0.1025557891954764
$ws.Range("T4").Value = 0.1025557891954764
